$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-14 from 2023-11-03 (45233)
# to 2023-11-13 (45243), keeping the existing date number format.
$ws.Range("C2:C14").Value = 45243
